$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 32, pushing the existing rows 32:67 down to 33:68.
$ws.Rows(32).Insert()

# Populate the newly inserted row 32 with the new weekly record.
$ws.Cells.Item(32, 1).Value = 9
$ws.Cells.Item(32, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(32, 3).Value = "Metropolitana"
$ws.Cells.Item(32, 4).Value = 44705
$ws.Cells.Item(32, 5).Value = 13
$ws.Cells.Item(32, 6).Value = "Fruta"
$ws.Cells.Item(32, 7).Value = 100102
$ws.Cells.Item(32, 8).Value = "Cítricos"
$ws.Cells.Item(32, 9).Value = 100102006
$ws.Cells.Item(32, 10).Value = "Pomelo"
$ws.Cells.Item(32, 11).Value = "Start Ruby"
$ws.Cells.Item(32, 12).Value = "Primera"
$ws.Cells.Item(32, 13).Value = 330
$ws.Cells.Item(32, 14).Value = 7500
$ws.Cells.Item(32, 15).Value = 7500
$ws.Cells.Item(32, 16).Value = 7500
$ws.Cells.Item(32, 17).Value = "`$/caja 14 kilos"
$ws.Cells.Item(32, 18).Value = "Región Metropolitana"
$ws.Cells.Item(32, 19).Value = 536
$ws.Cells.Item(32, 20).Value = 14
